$d = $word.ActiveDocument

# Locate the italic "(src/main/java/frc/team5333/recyclerush)" run by its
# exact original text so the edit is resilient to any surrounding changes.
$old = "(src/main/java/frc/team5333/recyclerush)"
$rng = $d.Content
$found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target run text"
}

$target = $d.Range($rng.Start, $rng.End)

# Split the single run into three runs, all keeping the original italic
# formatting, but with the path abbreviated to '...':
#   "(src/main/java" + "/..." + ")"
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:body>' +
       '<w:p>' +
       '<w:r w:rsidRPr="00C02D96"><w:rPr><w:i/></w:rPr><w:t>(src/main/java</w:t></w:r>' +
       '<w:r><w:rPr><w:i/></w:rPr><w:t>/...</w:t></w:r>' +
       '<w:r><w:rPr><w:i/></w:rPr><w:t>)</w:t></w:r>' +
       '</w:p>' +
       '</w:body>' +
       '</w:document>' +
       '</pkg:xmlData>' +
       '</pkg:part>' +
       '</pkg:package>'

$target.InsertXML($xml)
